# Insert a new data row at row 24, shifting existing rows 24:107 down to 25:108,
# then populate the new row 24 with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before current row 24 (pushes rows 24..107 down to 25..108)
$ws.Rows.Item(24).Insert()

# Populate the newly inserted row 24 with the new record
$ws.Range("A24").Value = 5
$ws.Range("B24").Value = "Macroferia Regional de Talca"
$ws.Range("C24").Value = "Maule"
$ws.Range("D24").Value = 44868
$ws.Range("E24").Value = 7
$ws.Range("F24").Value = 100112022
$ws.Range("G24").Value = "Arveja Verde"
$ws.Range("H24").Value = "Perfection"
$ws.Range("I24").Value = "Primera"
$ws.Range("J24").Value = 500
$ws.Range("K24").Value = 15000
$ws.Range("L24").Value = 15000
$ws.Range("M24").Value = 15000
$ws.Range("N24").Value = "`$/saco 25 kilos"
$ws.Range("O24").Value = "Región del Maule"
$ws.Range("P24").Value = 600
$ws.Range("Q24").Value = 25
$ws.Range("R24").Value = "Hortaliza"
